$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new lines (line7, line8) were added ahead of the "extr" entries, so the
# name column shifts: line7/line8 take the old extr1/extr2 rows, extr1..extr6
# shift up into extr1..extr8's old rows, and extr7/extr8 are brand new rows.
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10

$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8

$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# New rows 16 and 17
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

# Copy the "index" column styling (bold, centered, bordered) from the row above
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
